# Reproduce the "Add files via upload" commit against NovVacation.xlsx:
#  - bump a handful of transaction amounts
#  - fill in a vendor name on nov.eur and nov.mad
#  - leave cursor/selection state matching the final save
#  - switch the active sheet from xrates to nov.mad

$wb = $excel.ActiveWorkbook

# --- nov.gbp (sheet 1): amount 8.9 -> 10, cursor left on E3 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("E3").Value = 10

# --- nov.usd (sheet 2): amount 10.53 -> 10, cursor left on E3 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("E3").Value = 10

# --- nov.eur (sheet 3): vendor "Starbucks" added, amount 5.1 -> 100 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("C3").Value = "Starbucks"
$ws3.Range("E3").Value = 100

# --- nov.mad (sheet 4): vendor filled in, amount 0 -> 1000 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("C3").Value = "le meridien hotel"
$ws4.Range("E3").Value = 1000

# --- nov.inr (sheet 5): amount 1500 -> 1000, cursor left on E4 ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("E3").Value = 1000

# --- xrates (sheet 6): no longer the active tab ---
$ws6 = $wb.Worksheets.Item(6)

# Leave each sheet's selection where the author left it, then finish on
# nov.mad so it becomes the active/selected tab on save (taking that
# state over from xrates).
$ws1.Activate()
$ws1.Range("E3").Select()

$ws2.Activate()
$ws2.Range("E3").Select()

$ws3.Activate()
$ws3.Range("E3").Select()

$ws5.Activate()
$ws5.Range("E4").Select()

$ws6.Activate()
$ws6.Range("B13").Select()

$ws4.Activate()
$ws4.Range("A4").Select()
